$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template row to copy formatting (styles) from.
$templateRow = 101

$rowsData = @(
    @{
        Row = 102
        A = 101
        B = "portugal"
        C = "liga-portugal"
        D = "2023-2024"
        E = 45262.6875
        F = "Rio Ave"
        G = 1
        H = "Estrela"
        I = 1
        J = 2.26
        K = "15/11/2023 15:12"
        L = 2.06
        M = "02/12/2023 16:20"
        N = 3.36
        O = "15/11/2023 15:12"
        P = 3.43
        Q = "02/12/2023 16:20"
        R = 3.18
        S = "15/11/2023 15:12"
        T = 3.72
        U = "02/12/2023 16:20"
        V = "https://www.betexplorer.com/football/portugal/liga-portugal/rio-ave-estrela-da-amadora/vLYBGNdi/"
    },
    @{
        Row = 103
        A = 102
        B = "portugal"
        C = "liga-portugal"
        D = "2023-2024"
        E = 45262.6875
        F = "SC Farense"
        G = 1
        H = "Vitoria Guimaraes"
        I = 2
        J = 2.77
        K = "15/11/2023 15:12"
        L = 3.38
        M = "02/12/2023 16:28"
        N = 3.25
        O = "15/11/2023 15:12"
        P = 3.65
        Q = "02/12/2023 16:28"
        R = 2.59
        S = "15/11/2023 15:12"
        T = 2.1
        U = "02/12/2023 16:21"
        V = "https://www.betexplorer.com/football/portugal/liga-portugal/sc-farense-vitoria-guimaraes/8nofbtsN/"
    },
    @{
        Row = 104
        A = 103
        B = "portugal"
        C = "liga-portugal"
        D = "2023-2024"
        E = 45262.79166666666
        F = "Famalicao"
        G = 0
        H = "FC Porto"
        I = 3
        J = 5.83
        K = "15/11/2023 15:12"
        L = 7.45
        M = "02/12/2023 18:59"
        N = 4.21
        O = "15/11/2023 15:12"
        P = 4.46
        Q = "02/12/2023 18:56"
        R = 1.58
        S = "15/11/2023 15:12"
        T = 1.49
        U = "02/12/2023 18:56"
        V = "https://www.betexplorer.com/football/portugal/liga-portugal/famalicao-fc-porto/zZnn02CA/"
    }
)

foreach ($rd in $rowsData) {
    $r = $rd.Row

    # Copy the whole template row so that cell styles (column A index style,
    # column E date style, borders, bold, etc.) are carried over exactly.
    $ws.Range("A" + $templateRow + ":V" + $templateRow).Copy($ws.Range("A" + $r))

    $ws.Range("A" + $r).Value = $rd.A
    $ws.Range("B" + $r).Value = $rd.B
    $ws.Range("C" + $r).Value = $rd.C
    $ws.Range("D" + $r).Value = $rd.D
    $ws.Range("E" + $r).Value = $rd.E
    $ws.Range("F" + $r).Value = $rd.F
    $ws.Range("G" + $r).Value = $rd.G
    $ws.Range("H" + $r).Value = $rd.H
    $ws.Range("I" + $r).Value = $rd.I
    $ws.Range("J" + $r).Value = $rd.J
    $ws.Range("K" + $r).Value = $rd.K
    $ws.Range("L" + $r).Value = $rd.L
    $ws.Range("M" + $r).Value = $rd.M
    $ws.Range("N" + $r).Value = $rd.N
    $ws.Range("O" + $r).Value = $rd.O
    $ws.Range("P" + $r).Value = $rd.P
    $ws.Range("Q" + $r).Value = $rd.Q
    $ws.Range("R" + $r).Value = $rd.R
    $ws.Range("S" + $r).Value = $rd.S
    $ws.Range("T" + $r).Value = $rd.T
    $ws.Range("U" + $r).Value = $rd.U
    $ws.Range("V" + $r).Value = $rd.V
}

$excel.CutCopyMode = $false

Write-Host "Rows 102-104 added."
